$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Base" header labels (row 1, columns B:M) to the new
# plate/sample-prefixed naming scheme used for this PFAS experiment.
$ws.Range("B1").Value = "3920_Base"
$ws.Range("C1").Value = "3921_Base"
$ws.Range("D1").Value = "3922_Base"
$ws.Range("E1").Value = "3923_Base"
$ws.Range("F1").Value = "3932_Base"
$ws.Range("G1").Value = "3933_Base"
$ws.Range("H1").Value = "3934_Base"
$ws.Range("I1").Value = "3935_Base"
$ws.Range("J1").Value = "3936_Base"
$ws.Range("K1").Value = "3937_Base"
$ws.Range("L1").Value = "3938_Base"
$ws.Range("M1").Value = "3939_Base"

# Update "Final" header labels (row 1, columns Q:AB) likewise.
$ws.Range("Q1").Value = "3920_ Final"
$ws.Range("R1").Value = "3921_Final"
$ws.Range("S1").Value = "3922_Final"
$ws.Range("T1").Value = "3923_Final"
$ws.Range("U1").Value = "3932_Final"
$ws.Range("V1").Value = "3933_Final"
$ws.Range("W1").Value = "3934_Final"
$ws.Range("X1").Value = "3935_Final"
$ws.Range("Y1").Value = "3936_Final"
$ws.Range("Z1").Value = "3937_Final"
$ws.Range("AA1").Value = "3938_Final"
$ws.Range("AB1").Value = "3939_Final"

# Adjust the visible window/selection to match the saved view state.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("N7").Select()
